$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "DTSER-112"
$ws.Range("G2").Value = "NA"
$ws.Range("K2").Value = "El documento denominado “Ciénaga de Ayapel: Riqueza en Biodiversidad y Recursos Hídricos”, hace parte de la serie de Documentos de Trabajo sobre Economía Regional del Banco de la República, sucursal Cartagena. El objetivo principal de este documento es evaluar la riqueza en biodiversidad y los recursos hídricos de la ciénaga, así como identificar los desafíos ambientales y socioeconómicos que enfrenta la región."
$ws.Range("O2").Value = "No"
$ws.Range("P2").Value = "NA"
$ws.Range("Q2").Value = "El documento presenta un análisis exhaustivo , con datos relevantes sobre su biodiversidad, recursos hídricos y condiciones socioeconómicas. Sin embargo, algunos apartados podrían beneficiarse de una mayor profundidad en el análisis, especialmente en la sección de problemáticas ambientales, donde se mencionan desafíos sin ofrecer soluciones concretas o ejemplos de mejores prácticas. "
$ws.Range("R2").Value = "Ciénaga de Ayapel, biodiversidad, recursos hídricos, ecosistema, sostenibilidad, turismo ecológico, agricultura, ganadería, pesca artesanal, degradación ambiental, educación, condiciones sociales, conservación, desarrollo sostenible, Córdoba, Banco de la República, población, cultivos."

$ws.Range("H4").Select() | Out-Null
